$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 311-312, pushing the existing rows 311+ down to 313+.
$ws.Rows("311:312").Insert()

# New row 311: same market/product metadata as the (now shifted) row 313,
# but with an updated date and updated price figures.
$ws.Range("A311").Value = 5
$ws.Range("B311").Value = "Macroferia Regional de Talca"
$ws.Range("C311").Value = "Maule"
$ws.Range("D311").Value = 44461
$ws.Range("E311").Value = 7
$ws.Range("F311").Value = 100112020
$ws.Range("G311").Value = "Tomate"
$ws.Range("H311").Value = "Larga vida"
$ws.Range("I311").Value = "Primera"
$ws.Range("J311").Value = 1500
$ws.Range("K311").Value = 15000
$ws.Range("L311").Value = 15000
$ws.Range("M311").Value = 15000
$ws.Range("N311").Value = "$/bandeja 18 kilos"
$ws.Range("O311").Value = "Región de Arica y Parinacota"
$ws.Range("P311").Value = 833
$ws.Range("Q311").Value = 18
$ws.Range("R311").Value = "Hortaliza"

# New row 312.
$ws.Range("A312").Value = 5
$ws.Range("B312").Value = "Macroferia Regional de Talca"
$ws.Range("C312").Value = "Maule"
$ws.Range("D312").Value = 44461
$ws.Range("E312").Value = 7
$ws.Range("F312").Value = 100112020
$ws.Range("G312").Value = "Tomate"
$ws.Range("H312").Value = "Larga vida"
$ws.Range("I312").Value = "Primera"
$ws.Range("J312").Value = 1500
$ws.Range("K312").Value = 7000
$ws.Range("L312").Value = 7000
$ws.Range("M312").Value = 7000
$ws.Range("N312").Value = "$/caja 10 kilos"
$ws.Range("O312").Value = "Región de Arica y Parinacota"
$ws.Range("P312").Value = 700
$ws.Range("Q312").Value = 10
$ws.Range("R312").Value = "Hortaliza"
